$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixing Block Data")

# --- Set the values first (Piece ID "4" is forced to text because it is a
#     purely numeric-looking string, same as the rest of that column). ---

function Set-RowValues {
    param($Row, $Hole)

    $ws.Cells.Item($Row, 1).Value  = "2025-11-21 14:27:43"
    $ws.Cells.Item($Row, 2).Value  = "SA01"
    $ws.Cells.Item($Row, 3).Value  = "Mixing Block"
    $ws.Cells.Item($Row, 4).Value  = "A"

    $pieceCell = $ws.Cells.Item($Row, 5)
    $pieceCell.NumberFormat = "@"
    $pieceCell.Value = "4"

    $ws.Cells.Item($Row, 6).Value  = "IN"
    $ws.Cells.Item($Row, 7).Value  = $Hole
    $ws.Cells.Item($Row, 8).Value  = "Inner"
    $ws.Cells.Item($Row, 9).Value  = 111
    $ws.Cells.Item($Row, 10).Value = 4
    $ws.Cells.Item($Row, 11).Value = 3.5
    $ws.Cells.Item($Row, 12).Value = 4.5
    $ws.Cells.Item($Row, 13).Value = "FAIL"
}

Set-RowValues 29 "H1"
Set-RowValues 30 "H2"

# --- Apply formatting by cloning the existing FAIL-row style (row 28) so the
#     workbook's existing fill/border/style entries are reused as-is. ---

$srcFmt = $ws.Range("A28:O28")
$row29 = $ws.Range("A29:O29")
$row30 = $ws.Range("A30:O30")

$srcFmt.Copy()
$row29.PasteSpecial(-4122)
$row30.PasteSpecial(-4122)

# Row 29 is not the last row of this submission batch, so (matching the rest
# of the sheet's convention) it keeps the FAIL fill but loses the bottom
# border that marks the end of a batch; row 30 keeps it.
$row29.Borders.Item(9).LineStyle = -4142

$excel.CutCopyMode = 0
